# Updated cryptos list -- refresh Price (column D) and Volume(1h) (column E)
# figures for the coinranking.com snapshot on the "cryptos" sheet.
#
# Column D (Price) is stored as literal text in the workbook, not as a
# number -- every existing D-cell is t="inlineStr"/shared-string text,
# including values that look like plain numbers (e.g. "211.82"). A bare
# .Value assignment of a numeric-looking string lets Excel reinterpret it
# as a Number (dropping significant trailing zeros, e.g. "1.00" -> 1), so
# every Price write below is prefixed with a single quote -- exactly as
# typing '211.71 into a cell would -- to force Text storage while keeping
# the stored/displayed string identical to the target value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — Bitcoin
$ws.Range("D2").Value = "'26.715.15"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3 — Ethereum
$ws.Range("D3").Value = "'1.601.05"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4 — TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 — BNB
$ws.Range("D5").Value = "'211.71"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6 — XRP
$ws.Range("E6").Value = "  -0.45%  "

# Row 7 — USDC
$ws.Range("E7").Value = "  +0.15%  "

# Row 8 — Dogecoin
$ws.Range("D8").Value = "'0.0619"
$ws.Range("E8").Value = "  +0.22%  "

# Row 9 — Cardano
$ws.Range("E9").Value = "  +0.18%  "

# Row 10 — Solana
$ws.Range("D10").Value = "'19.70"
$ws.Range("E10").Value = "  +0.87%  "

# Row 11 — TRON
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  +1.02%  "

# Row 12 — Wrapped liquid staked Ether 2.0
$ws.Range("D12").Value = "'1.827.23"
$ws.Range("E12").Value = "  +0.31%  "

# Row 13 — Wrapped Ether
$ws.Range("D13").Value = "'1.601.32"
$ws.Range("E13").Value = "  +0.37%  "

# Row 14 — Polkadot
$ws.Range("D14").Value = "'4.04"
$ws.Range("E14").Value = "  +0.46%  "

# Row 15 — Polygon
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  -0.24%  "

# Row 16 — Litecoin
$ws.Range("D16").Value = "'65.06"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17 — Shiba Inu
$ws.Range("D17").Value = "'0.0₃0739"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18 — Bitcoin Cash
$ws.Range("D18").Value = "'210.03"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19 — Dai
$ws.Range("D19").Value = "'1.00"
$ws.Range("E19").Value = "  +0.13%  "

# Row 20 — Chainlink
$ws.Range("D20").Value = "'7.17"
$ws.Range("E20").Value = "  +2.48%  "

# Row 21 — Uniswap
$ws.Range("E21").Value = "  -0.21%  "

# Row 22 — Toncoin
$ws.Range("E22").Value = "  -2.69%  "

# Row 23 — Avalanche
$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24 — Monero
$ws.Range("D24").Value = "'143.62"

# Row 25 — BinanceUSD
$ws.Range("E25").Value = "  -0.11%  "

# Row 26 — Cosmos
$ws.Range("D26").Value = "'7.09"
$ws.Range("E26").Value = "  -0.13%  "

# Row 27 — Stellar
$ws.Range("E27").Value = "  -0.94%  "

# Row 28 — Ethereum Classic
$ws.Range("D28").Value = "'15.37"
$ws.Range("E28").Value = "  +0.58%  "

# Row 29 — Hedera
$ws.Range("E29").Value = "  -0.61%  "

# Row 30 — PancakeSwap
$ws.Range("E30").Value = "  +0.03%  "

# Row 31 — Filecoin
$ws.Range("E31").Value = "  +1.13%  "

# Row 33 — Maker
$ws.Range("D33").Value = "'1.292.34"
$ws.Range("E33").Value = "  +0.84%  "

# Row 34 — Huobi Token
$ws.Range("E34").Value = "  +0.57%  "

# Row 35 — Lido DAO Token
$ws.Range("E35").Value = "  +0.49%  "

# Row 36 — Immutable X
$ws.Range("D36").Value = "'0.604"
$ws.Range("E36").Value = "  -2.27%  "

# Row 37 — WEMIX Token
$ws.Range("E37").Value = "  +9.44%  "

# Row 38 — VeChain
$ws.Range("E38").Value = "  -0.15%  "

# Row 39 — ARBITRUM
$ws.Range("D39").Value = "'0.831"
$ws.Range("E39").Value = "  -0.33%  "

# Row 40 — Frax Share
$ws.Range("E40").Value = "  -2.28%  "

# Row 41 — MX Token
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 — Trust Wallet Token
$ws.Range("D42").Value = "'0.786"
$ws.Range("E42").Value = "  +0.10%  "

# Row 43 — Aave
$ws.Range("D43").Value = "'62.65"
$ws.Range("E43").Value = "  -1.98%  "

# Row 44 — Rocket Pool ETH
$ws.Range("D44").Value = "'1.738.65"
$ws.Range("E44").Value = "  +0.24%  "

# Row 45 — Quant
$ws.Range("D45").Value = "'90.60"
$ws.Range("E45").Value = "  -0.29%  "

# Row 47 — Algorand
$ws.Range("D47").Value = "'0.101"
$ws.Range("E47").Value = "  -0.23%  "

# Row 48 — Cronos
$ws.Range("E48").Value = "  +1.69%  "

# Row 49 — USDD
$ws.Range("E49").Value = "  +0.20%  "

# Row 50 — EnergySwap
$ws.Range("D50").Value = "'7.40"
$ws.Range("E50").Value = "  +0.33%  "

# Row 51 — Mantle
$ws.Range("E51").Value = "  +0.83%  "
